$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 235-237: the three matches played on 2024-03-08 get re-ordered
# (cyclic rotation of rows 235 -> 236 -> 237 -> 235), while the leading "id"
# column (A) keeps counting 233,234,235 in place.
# ---------------------------------------------------------------------------

# Row 235 (id stays 233)
$ws.Range("B235").Value = 6852370
$ws.Range("F235").Value = "Dinamo Bucharest"
$ws.Range("G235").Value = "ACS UTA Batrana Doamna"
$ws.Range("I235").Value = 0
$ws.Range("J235").Value = "H"
$ws.Range("K235").Value = 2.55
$ws.Range("L235").Value = 2.875
$ws.Range("M235").Value = 3
$ws.Range("N235").Value = 2.375
$ws.Range("O235").Value = 3
$ws.Range("P235").Value = 3.1
$ws.Range("Q235").Value = -0.25
$ws.Range("R235").Value = 2
$ws.Range("S235").Value = 1.85
$ws.Range("T235").Value = 2.25
$ws.Range("U235").Value = 1.975
$ws.Range("V235").Value = 1.875
$ws.Range("W235").Value = 1.375
$ws.Range("X235").Value = -1
$ws.Range("Y235").Value = -1
$ws.Range("Z235").Value = 1
$ws.Range("AA235").Value = -1
$ws.Range("AB235").Value = -1
$ws.Range("AC235").Value = 0.875

# Row 236 (id stays 234)
$ws.Range("B236").Value = 6836277
$ws.Range("F236").Value = "CFR Cluj"
$ws.Range("G236").Value = "AFC Hermannstadt"
$ws.Range("I236").Value = 0
$ws.Range("J236").Value = "H"
$ws.Range("K236").Value = 1.7
$ws.Range("L236").Value = 3.4
$ws.Range("M236").Value = 5
$ws.Range("N236").Value = 1.65
$ws.Range("O236").Value = 3.5
$ws.Range("P236").Value = 5.25
$ws.Range("Q236").Value = -0.75
$ws.Range("R236").Value = 1.85
$ws.Range("S236").Value = 2
$ws.Range("T236").Value = 2.25
$ws.Range("U236").Value = 1.875
$ws.Range("V236").Value = 1.975
$ws.Range("W236").Value = 0.6499999999999999
$ws.Range("X236").Value = -1
$ws.Range("Y236").Value = -1
$ws.Range("Z236").Value = 0.425
$ws.Range("AA236").Value = -0.5
$ws.Range("AB236").Value = -1
$ws.Range("AC236").Value = 0.9750000000000001

# Row 237 (id stays 235)
$ws.Range("B237").Value = 6870268
$ws.Range("F237").Value = "Petrolul Ploiesti"
$ws.Range("G237").Value = "ACS Sepsi"
$ws.Range("I237").Value = 2
$ws.Range("J237").Value = "A"
$ws.Range("K237").Value = 2.8
$ws.Range("L237").Value = 3
$ws.Range("M237").Value = 2.55
$ws.Range("N237").Value = 3
$ws.Range("O237").Value = 3.2
$ws.Range("P237").Value = 2.3
$ws.Range("Q237").Value = 0.25
$ws.Range("R237").Value = 1.85
$ws.Range("S237").Value = 2
$ws.Range("T237").Value = 2.25
$ws.Range("U237").Value = 1.875
$ws.Range("V237").Value = 1.975
$ws.Range("W237").Value = -1
$ws.Range("X237").Value = -1
$ws.Range("Y237").Value = 1.3
$ws.Range("Z237").Value = -1
$ws.Range("AA237").Value = 1
$ws.Range("AB237").Value = 0.875
$ws.Range("AC237").Value = -1

# ---------------------------------------------------------------------------
# Rows 251-253: refresh of the upcoming-fixtures list. The matches that were
# in rows 251-253 (ids 8010912 / 8010913 / 7951749) are gone, rows 254-255
# shift up into 251-252, and a brand-new fixture is appended as row 253.
# ---------------------------------------------------------------------------

# Row 251 (id stays 249)
$ws.Range("B251").Value = 7951779
$ws.Range("E251").Value = 45382.33333333334
$ws.Range("F251").Value = "FC U Craiova 1948"
$ws.Range("G251").Value = "Otelul Galati"
$ws.Range("K251").Value = 2.3
$ws.Range("L251").Value = 3.2
$ws.Range("M251").Value = 3.2
$ws.Range("N251").Value = 2.3
$ws.Range("O251").Value = 3.2
$ws.Range("P251").Value = 3.2
$ws.Range("Q251").Value = -0.25
$ws.Range("R251").Value = 2
$ws.Range("S251").Value = 1.85
$ws.Range("T251").Value = 2.25
$ws.Range("U251").Value = 2.1
$ws.Range("V251").Value = 1.775

# Row 252 (id stays 250)
$ws.Range("B252").Value = 7951748
$ws.Range("E252").Value = 45382.625
$ws.Range("F252").Value = "Farul Constanta"
$ws.Range("G252").Value = "FCSB"
$ws.Range("K252").Value = 3.6
$ws.Range("L252").Value = 3.3
$ws.Range("M252").Value = 2
$ws.Range("N252").Value = 3.6
$ws.Range("O252").Value = 3.3
$ws.Range("P252").Value = 2.05
$ws.Range("Q252").Value = 0.25
$ws.Range("R252").Value = 2.025
$ws.Range("S252").Value = 1.825
$ws.Range("T252").Value = 2.25
$ws.Range("U252").Value = 1.8
$ws.Range("V252").Value = 2.05

# Row 253 (id stays 251) - brand new fixture replacing the old one
$ws.Range("B253").Value = 7951780
$ws.Range("E253").Value = 45383.60416666666
$ws.Range("F253").Value = "Dinamo Bucharest"
$ws.Range("G253").Value = "Petrolul Ploiesti"
$ws.Range("K253").Value = 2.3
$ws.Range("L253").Value = 3
$ws.Range("M253").Value = 3.4
$ws.Range("N253").Value = 2.375
$ws.Range("O253").Value = 3
$ws.Range("P253").Value = 3.3
$ws.Range("Q253").Value = -0.25
$ws.Range("R253").Value = 2.05
$ws.Range("S253").Value = 1.8
$ws.Range("T253").Value = 2
$ws.Range("U253").Value = 1.95
$ws.Range("V253").Value = 1.9

# ---------------------------------------------------------------------------
# Remove the now-obsolete trailing rows 254-255 (their data moved up into
# rows 251-252 above, and row 253 got fresh data), shrinking the sheet's
# used range down to A1:AC253.
# ---------------------------------------------------------------------------
$ws.Rows("254:255").Delete()
